$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the bold/bordered
# header style used by the existing headers (style index 1).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data columns I and J for rows 2-29.
$data = @(
    @(2, 7, 7),
    @(3, 8, 8),
    @(4, 7, 7),
    @(5, 10, 10),
    @(6, 7, 8),
    @(7, 9, 10),
    @(8, 9, 9),
    @(9, 7, 7),
    @(10, 8, 9),
    @(11, 8, 9),
    @(12, 8, 8),
    @(13, 8, 8),
    @(14, 7, 8),
    @(15, 7, 7),
    @(16, 7, 7),
    @(17, 12, 12),
    @(18, 7, 7),
    @(19, 9, 9),
    @(20, 9, 9),
    @(21, 7, 7),
    @(22, 7, 7),
    @(23, 6, 6),
    @(24, 7, 7),
    @(25, 7, 7),
    @(26, 7, 7),
    @(27, 6, 6),
    @(28, 6, 6),
    @(29, 3, 3)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 9).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
}
